$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D7").Value = -7.713000000000001
$ws.Range("A9").Value = -21.778
$ws.Range("D12").Value = -7.269
$ws.Range("D14").Value = -7.583
$ws.Range("A18").Value = -22.387
$ws.Range("A20").Value = -20.1
$ws.Range("D26").Value = -8.122999999999999
$ws.Range("A27").Value = -21.86
$ws.Range("D27").Value = -8.599
$ws.Range("D29").Value = -7.292
$ws.Range("A35").Value = -19.823
$ws.Range("D37").Value = -7.822
$ws.Range("D38").Value = -7.722
$ws.Range("D51").Value = -8.638000000000002
$ws.Range("D52").Value = -8.1
$ws.Range("D55").Value = -7.918000000000001
$ws.Range("A69").Value = -21.627
$ws.Range("D69").Value = -7.220999999999999
$ws.Range("D70").Value = -7.269
$ws.Range("A76").Value = -20.047
$ws.Range("A78").Value = -20.364
$ws.Range("D81").Value = -7.622
$ws.Range("A82").Value = -22.095
$ws.Range("A83").Value = -21.802
$ws.Range("D83").Value = -8.462999999999999
$ws.Range("A93").Value = -21.524
$ws.Range("D102").Value = -7.906999999999999
